# 10Th - MB for single stock and added new group
# This adds two new weekly snapshot columns (Jun_27, Jun_26) ahead of the
# existing Jun_13 / Jun_10 columns, records a new analyst rating change
# (ValuEngine: Buy -> Hold on 6/21/2018) and appends two new analyst rows
# (Benchmark, Evercore ISI) to the bottom of the tracker.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room for the new weekly columns ----------------------------
# Existing layout is A (analyst), B (Jun_13), C (Jun_10).
# Insert three new columns before B so the old B/C shift right to E/F.
$ws.Columns("B:D").Insert()

# Match the look & feel of the neighbouring "Jun_13" column for the new ones.
$ws.Columns("B").ColumnWidth = 29.1640625
$ws.Columns("C").ColumnWidth = 29.1640625
$ws.Columns("D").ColumnWidth = 29.1640625

# --- 2. New week headers -------------------------------------------------
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# --- 3. Default every analyst row to "UN" (unrated) for the new weeks ----
$analysts = @{
    2  = "Needham & Company LLC"
    3  = "Credit Suisse Group"
    4  = "Morningstar"
    5  = "Zacks Investment Research"
    6  = "Goldman Sachs Group"
    7  = "Morgan Stanley"
    8  = "KeyCorp"
    9  = "Argus"
    10 = "Piper Jaffray Companies"
    11 = "ValuEngine"
    12 = "Sidoti"
    13 = "Jefferies Financial Group"
    14 = "Stifel Nicolaus"
    15 = "Vetr"
    16 = "Fundamental Research"
    17 = "JPMorgan Chase & Co."
    18 = "Royal Bank of Canada"
    19 = "Citigroup"
    20 = "Bank of America"
    21 = "Wells Fargo & Co"
    22 = "BidaskClub"
    23 = "Roth Capital"
    24 = "Janney Montgomery Scott"
    25 = "William Blair"
    26 = "Stephens"
    27 = "Barclays"
}

for ($row = 2; $row -le 27; $row++) {
    $ws.Range("B$row").Value = "UN"
    $ws.Range("C$row").Value = "UN"
    $ws.Range("D$row").Value = "UN"
}

# --- 4. Record the new rating change (ValuEngine downgraded Buy -> Hold) --
$change = "6/21/2018,Downgrades,Buy -> Hold,"
$ws.Range("C11").Value = $change
$ws.Range("D11").Value = $change
$ws.Range("D11").Interior.ColorIndex = 45

# --- 5. Append the two newly tracked analysts -----------------------------
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"
$ws.Range("C28").Value = "UN"
$ws.Range("D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
$ws.Range("C29").Value = "UN"
$ws.Range("D29").Value = "UN"

Write-Output "edit complete"
